$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new daily row (row 82) with the next day's values
$newRow = 82

$ws.Cells.Item($newRow, 1).Value = 45669
$ws.Cells.Item($newRow, 2).Value = 192
$ws.Cells.Item($newRow, 3).Value = 190
$ws.Cells.Item($newRow, 4).Value = 191

# Move the "last row" date style (currently on A81) onto the new last row (A82),
# and reset A81 back to the normal date style used by all other rows (same as A80).
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item(81, 1).NumberFormat
$ws.Cells.Item(81, 1).NumberFormat = $ws.Cells.Item(80, 1).NumberFormat
